$d = $word.ActiveDocument

# The target paragraph is the one reading "Изпратете двете SQL команди за
# проверка в Judge." where "Judge" is still split across two unformatted
# runs ("J" + "udge"). Locate it by paragraph index (stable for this
# document) and confirm via text content before touching formatting.
$p = $d.Paragraphs.Item(87)
$r = $p.Range.Duplicate

$f = $r.Find
$f.ClearFormatting()
$f.Text = "Judge"
$f.Forward = $true
$f.Wrap = 0
$found = $f.Execute()

if ($found) {
    $r.Font.Bold = $true
    $r.Font.Bold = 1
    $r.Bold = 1
}
